$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Purchase 22-23"): row 17 gets replaced by what used to be
# row 19's data (Sr. No in column A is untouched), and row 19 is removed.

# Bring G19's number format/style onto the not-yet-existing G17 first, so
# the value we set below lands with the right style (matches G19's style).
$ws1.Cells.Item(19, 7).Copy()
$ws1.Cells.Item(17, 7).PasteSpecial(-4122)

$ws1.Cells.Item(17, 2).Value = $ws1.Cells.Item(19, 2).Value()
$ws1.Cells.Item(17, 3).ClearContents()
$ws1.Cells.Item(17, 4).Value = $ws1.Cells.Item(19, 4).Value()
$ws1.Cells.Item(17, 5).Value = $ws1.Cells.Item(19, 5).Value()
$ws1.Cells.Item(17, 7).Value = $ws1.Cells.Item(19, 7).Value()

# Row 17's F column formula already reads "=E17" so it recalculates on
# its own once E17 is updated above; row 19 can now be dropped entirely.
$ws1.Rows.Item(19).Delete()

# --- Selection / active-sheet bookkeeping -----------------------------
# Sheet2 ("Sale 22-23") keeps a lingering selection but is no longer the
# active tab, so update its selection first (selecting a range on a sheet
# also activates that sheet) and finish on Sheet1 so it ends up active.
$ws2.Range("B31").Select()
$ws1.Range("A18").Select()
